# Updated cryptos list values (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking price strings
# (e.g. "1.004") are stored as literal text, not coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.409.07"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.858.61"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "314.90"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "0.4634"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.07325"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "0.8898"
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "0.07831"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "1.894.50"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "5.413"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "6.574"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "91.97"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "0.000008991"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "27.418.85"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "5.136"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").Value = "2.124.09"
$ws.Range("E24").Value = "  +5.69%  "
$ws.Range("D25").Value = "1.939"
$ws.Range("E25").Value = "  +5.46%  "
$ws.Range("D26").Value = "152.30"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "18.48"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "116.34"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("D31").Value = "0.08855"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "3.121"
$ws.Range("E32").Value = "  +5.14%  "
$ws.Range("D33").Value = "0.7692"
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("E34").Value = "  +3.32%  "
$ws.Range("D35").Value = "4.516"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").Value = "2.690"
$ws.Range("E36").Value = "  +8.58%  "
$ws.Range("D37").Value = "1.083"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").Value = "0.05249"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "2.966"
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").Value = "7.087"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("D42").Value = "0.5158"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "0.1642"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "8.409"
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("D45").Value = "0.4822"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "10.32"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "103.18"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "1.656"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "65.40"
$ws.Range("E51").Value = "  +0.51%  "

# Restore the original (default) style on column D so no stray
# number-format style reference is left on the cells.
$ws.Range("D2:D51").Style = "Normal"

